$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.67"
$ws.Range("E2").Value = "'3.17%"
$ws.Range("D3").Value = "'41.55"
$ws.Range("E3").Value = "'2.77%"
$ws.Range("D4").Value = "'5.012"
$ws.Range("E4").Value = "'-0.24%"
$ws.Range("D5").Value = "'0.07542"
$ws.Range("E5").Value = "'3.48%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.365"
$ws.Range("E6").Value = "'1.90%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.572"
$ws.Range("E7").Value = "'3.29%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9301"
$ws.Range("E8").Value = "'1.39%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.401"
$ws.Range("E9").Value = "'0.17%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1205"
$ws.Range("E10").Value = "'2.33%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1839"
$ws.Range("E11").Value = "'6.60%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08844"
$ws.Range("E12").Value = "'2.24%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04078"
$ws.Range("E13").Value = "'-2.38%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1054"
$ws.Range("E14").Value = "'0.08%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001273"
$ws.Range("E15").Value = "'0.48%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005840"
$ws.Range("E16").Value = "'1.05%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.343"
$ws.Range("E17").Value = "'-1.38%"
$ws.Range("D18").Value = "'0.3336"
$ws.Range("E18").Value = "'0.52%"
$ws.Range("D19").Value = "'7.926"
$ws.Range("E19").Value = "'0.55%"
$ws.Range("E20").Value = "'5.54%"
$ws.Range("D21").Value = "'0.2994"
$ws.Range("E21").Value = "'3.71%"
$ws.Range("D22").Value = "'0.04041"
$ws.Range("E22").Value = "'4.34%"
$ws.Range("D23").Value = "'0.001262"
$ws.Range("E23").Value = "'-0.63%"
$ws.Range("D24").Value = "'0.003907"
$ws.Range("E24").Value = "'1.39%"
$ws.Range("D25").Value = "'0.0001228"
$ws.Range("E25").Value = "'-4.26%"
$ws.Range("E26").Value = "'-0.14%"
$ws.Range("D38").Value = "'0.02424"
$ws.Range("E38").Value = "'5.01%"
$ws.Range("D39").Value = "'0.05218"
$ws.Range("E39").Value = "'5.32%"
$ws.Range("D40").Value = "'0.005879"
$ws.Range("E40").Value = "'-12.81%"
$ws.Range("D41").Value = "'0.007801"
$ws.Range("E41").Value = "'1.63%"
$ws.Range("D42").Value = "'0.1332"
$ws.Range("E42").Value = "'4.68%"
$ws.Range("D43").Value = "'0.007354"
$ws.Range("E43").Value = "'-0.05%"
$ws.Range("D44").Value = "'0.007820"
$ws.Range("E44").Value = "'10.61%"
$ws.Range("D45").Value = "'0.2976"
$ws.Range("E45").Value = "'-4.39%"
$ws.Range("D46").Value = "'0.00006322"
$ws.Range("E46").Value = "'-1.82%"
$ws.Range("E47").Value = "'-0.36%"
$ws.Range("D48").Value = "'0.04510"
$ws.Range("E48").Value = "'28.86%"
$ws.Range("E50").Value = "'-0.36%"
$ws.Range("E51").Value = "'-0.36%"
